$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update 想去人数 (column F) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 639
$wsExhibit.Range("F3").Value = 2230
$wsExhibit.Range("F5").Value = 13402
$wsExhibit.Range("F11").Value = 1008
$wsExhibit.Range("F12").Value = 13829
$wsExhibit.Range("F21").Value = 11
$wsExhibit.Range("F22").Value = 46
$wsExhibit.Range("F23").Value = 5
$wsExhibit.Range("F24").Value = 1111
$wsExhibit.Range("F25").Value = 118
$wsExhibit.Range("F27").Value = 5558
$wsExhibit.Range("F31").Value = 32
$wsExhibit.Range("F32").Value = 23
$wsExhibit.Range("F33").Value = 142

# Sheet "全部类型" (fourth sheet) - update 想去人数 (column F) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 639
$wsAll.Range("F3").Value = 2230
$wsAll.Range("F5").Value = 13402
$wsAll.Range("F12").Value = 1008
$wsAll.Range("F13").Value = 13829
$wsAll.Range("F22").Value = 11
$wsAll.Range("F23").Value = 46
$wsAll.Range("F24").Value = 5
$wsAll.Range("F25").Value = 1111
$wsAll.Range("F26").Value = 118
$wsAll.Range("F28").Value = 5558
$wsAll.Range("F32").Value = 32
$wsAll.Range("F33").Value = 23
$wsAll.Range("F34").Value = 142
